$wb = $excel.ActiveWorkbook

# --- Insert 5 new sheets (Confidence, Joy, Anger, Surprise, Sorrow) between
# --- "BPP_win" and "Meta", merging in results from a previous emotion study.
$bppWin = $wb.Worksheets.Item("BPP_win")

$names = @("Confidence", "Joy", "Anger", "Surprise", "Sorrow")
$prev = $bppWin
foreach ($n in $names) {
    $newSheet = $wb.Worksheets.Add($null, $prev)
    $newSheet.Name = $n
    $prev = $newSheet
}

# --- Confidence: overall win/lose confidence prior ---
$ws = $wb.Worksheets.Item("Confidence")
$ws.Cells.Item(1,1).Value = 0
$ws.Cells.Item(1,2).Value = 1
$ws.Cells.Item(2,1).Value = 0.5
$ws.Cells.Item(2,2).Value = 0.5

# --- Joy: P(emotion level 1-5 | win=0/1) ---
$ws = $wb.Worksheets.Item("Joy")
$ws.Cells.Item(1,1).Value = "Confidence"
$ws.Cells.Item(1,2).Value = 1
$ws.Cells.Item(1,3).Value = 2
$ws.Cells.Item(1,4).Value = 3
$ws.Cells.Item(1,5).Value = 4
$ws.Cells.Item(1,6).Value = 5
$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,2).Value = 0.8
$ws.Cells.Item(2,3).Value = 0
$ws.Cells.Item(2,4).Value = 0
$ws.Cells.Item(2,5).Value = 0.1
$ws.Cells.Item(2,6).Value = 0.1
$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(3,2).Value = 0.1
$ws.Cells.Item(3,3).Value = 0.1
$ws.Cells.Item(3,4).Value = 0.1
$ws.Cells.Item(3,5).Value = 0
$ws.Cells.Item(3,6).Value = 0.7

# --- Anger: P(emotion level 1-5 | win=0/1) ---
$ws = $wb.Worksheets.Item("Anger")
$ws.Cells.Item(1,1).Value = "Confidence"
$ws.Cells.Item(1,2).Value = 1
$ws.Cells.Item(1,3).Value = 2
$ws.Cells.Item(1,4).Value = 3
$ws.Cells.Item(1,5).Value = 4
$ws.Cells.Item(1,6).Value = 5
$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,2).Value = 0.6
$ws.Cells.Item(2,3).Value = 0.2
$ws.Cells.Item(2,4).Value = 0
$ws.Cells.Item(2,5).Value = 0.2
$ws.Cells.Item(2,6).Value = 0
$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(3,2).Value = 0.9
$ws.Cells.Item(3,3).Value = 0
$ws.Cells.Item(3,4).Value = 0.1
$ws.Cells.Item(3,5).Value = 0
$ws.Cells.Item(3,6).Value = 0

# --- Surprise: P(emotion level 1-5 | win=0/1) ---
$ws = $wb.Worksheets.Item("Surprise")
$ws.Cells.Item(1,1).Value = "Confidence"
$ws.Cells.Item(1,2).Value = 1
$ws.Cells.Item(1,3).Value = 2
$ws.Cells.Item(1,4).Value = 3
$ws.Cells.Item(1,5).Value = 4
$ws.Cells.Item(1,6).Value = 5
$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,2).Value = 0.4
$ws.Cells.Item(2,3).Value = 0.2
$ws.Cells.Item(2,4).Value = 0.2
$ws.Cells.Item(2,5).Value = 0.1
$ws.Cells.Item(2,6).Value = 0.1
$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(3,2).Value = 0.5
$ws.Cells.Item(3,3).Value = 0.3
$ws.Cells.Item(3,4).Value = 0.1
$ws.Cells.Item(3,5).Value = 0.1
$ws.Cells.Item(3,6).Value = 0

# --- Sorrow: P(emotion level 1-5 | win=0/1) ---
$ws = $wb.Worksheets.Item("Sorrow")
$ws.Cells.Item(1,1).Value = "Confidence"
$ws.Cells.Item(1,2).Value = 1
$ws.Cells.Item(1,3).Value = 2
$ws.Cells.Item(1,4).Value = 3
$ws.Cells.Item(1,5).Value = 4
$ws.Cells.Item(1,6).Value = 5
$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,2).Value = 0.7
$ws.Cells.Item(2,3).Value = 0.1
$ws.Cells.Item(2,4).Value = 0.1
$ws.Cells.Item(2,5).Value = 0.1
$ws.Cells.Item(2,6).Value = 0
$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(3,2).Value = 1
$ws.Cells.Item(3,3).Value = 0
$ws.Cells.Item(3,4).Value = 0
$ws.Cells.Item(3,5).Value = 0
$ws.Cells.Item(3,6).Value = 0

# --- Restore per-sheet selections (matches the saved workbook state) ---
# Order matters: the sheet selected last becomes the active tab, and
# "Sorrow" is the author's final active tab, so it must be selected last.
[void]$wb.Worksheets.Item("Confidence").Range("A2").Select()
[void]$wb.Worksheets.Item("Joy").Range("G1").Select()
[void]$wb.Worksheets.Item("Anger").Range("G1").Select()
[void]$wb.Worksheets.Item("Surprise").Range("F1").Select()
[void]$wb.Worksheets.Item("Sorrow").Range("G2").Select()
